$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.4609249367666753, 0.7037408350246651, 0.6939948720570934, 0.833063546229874, 0.7182889450858672, 15),
    @(0.3216291049761831, 0.6088832333399757, 0.584052044952158, 0.7642329781893464, 0.7194279080246428, 14),
    @(0.2154567245497249, 0.5333824792787472, 0.4853342657017765, 0.6966593613106599, 0.6895569019380638, 13),
    @(0.4452680593528981, 0.5566372887181056, 0.4832252965807901, 0.6951440833243063, 0.5575547117719994, 12),
    @(0.3156011815745178, 0.4418649605082216, 0.2461552685847581, 0.4961403718553431, 0.4015050174666302, 11),
    @(0.3772995006961687, 0.4404536501459158, 0.2519683207200895, 0.5019644616106697, 0.348988072853359, 10),
    @(0.3543917895297701, 0.4159301033722839, 0.2171525390636058, 0.4659962865341373, 0.3209421651944155, 9),
    @(0.380213756763938, 0.4480587750649256, 0.2483796526124562, 0.4983770185436486, 0.3444534416474213, 8),
    @(0.3348409171772554, 0.4048669237504626, 0.2023859592503501, 0.4498732702110118, 0.3245182676007908, 7),
    @(0.3934609991409326, 0.465587436599884, 0.2489390515396543, 0.4989379235332331, 0.3360847994682415, 6)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}
